# Increment the "want to go" count (column F) by 1 for a handful of events.
# The same events show up both in their category sheet "展览" and in the
# aggregated "全部类型" sheet, so each needs to be updated in both places.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "展览";     Cell = "F16"; Value = 36 },
    @{ Sheet = "展览";     Cell = "F19"; Value = 1326 },
    @{ Sheet = "展览";     Cell = "F24"; Value = 7219 },
    @{ Sheet = "展览";     Cell = "F27"; Value = 36 },
    @{ Sheet = "展览";     Cell = "F29"; Value = 11 },
    @{ Sheet = "展览";     Cell = "F30"; Value = 5858 },
    @{ Sheet = "全部类型"; Cell = "F16"; Value = 36 },
    @{ Sheet = "全部类型"; Cell = "F19"; Value = 1326 },
    @{ Sheet = "全部类型"; Cell = "F25"; Value = 7219 },
    @{ Sheet = "全部类型"; Cell = "F28"; Value = 36 },
    @{ Sheet = "全部类型"; Cell = "F30"; Value = 11 },
    @{ Sheet = "全部类型"; Cell = "F32"; Value = 5858 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
